# Insert a new data row at row 28 (pushing existing rows 28..112 down to 29..113)
# and populate it with the new "Ecuador" price observation for the date 44715
# (2022-06-03). This mirrors the source diff which shows the whole data block
# from row 28 onward shifting down by one row, with a brand new record
# occupying the (now vacant) row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(28).Insert()

$ws.Cells.Item(28, 1).Value  = 11
$ws.Cells.Item(28, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value  = "Bíobío"
$ws.Cells.Item(28, 4).Value  = 44715
$ws.Cells.Item(28, 5).Value  = 8
$ws.Cells.Item(28, 6).Value  = "Fruta"
$ws.Cells.Item(28, 7).Value  = 100108
$ws.Cells.Item(28, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(28, 9).Value  = 100108002
$ws.Cells.Item(28, 10).Value = "Mango"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 300
$ws.Cells.Item(28, 14).Value = 9000
$ws.Cells.Item(28, 15).Value = 10000
$ws.Cells.Item(28, 16).Value = 9333
$ws.Cells.Item(28, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(28, 18).Value = "Ecuador"
$ws.Cells.Item(28, 19).Value = 2333
$ws.Cells.Item(28, 20).Value = 4
